$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.756.16'
$ws.Range("E2").Value = '  +4.31%  '
$ws.Range("D3").Value = '2.266.35'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '303.79'
$ws.Range("E5").Value = '  +3.11%  '
$ws.Range("D6").Value = '91.35'
$ws.Range("E6").Value = '  +4.67%  '
$ws.Range("E7").Value = '  +3.49%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  +2.39%  '
$ws.Range("D10").Value = '32.23'
$ws.Range("E10").Value = '  +4.84%  '
$ws.Range("D11").Value = '53.43'
$ws.Range("E12").Value = '  +1.86%  '
$ws.Range("E13").Value = '  +1.60%  '
$ws.Range("D14").Value = '6.57'
$ws.Range("E14").Value = '  +3.02%  '
$ws.Range("D15").Value = '2.617.41'
$ws.Range("E15").Value = '  +2.36%  '
$ws.Range("D16").Value = '14.16'
$ws.Range("E16").Value = '  +2.56%  '
$ws.Range("D17").Value = '2.258.82'
$ws.Range("E17").Value = '  +0.27%  '
$ws.Range("E18").Value = '  +3.50%  '
$ws.Range("D19").Value = '41.690.79'
$ws.Range("E19").Value = '  +4.31%  '
$ws.Range("D20").Value = '12.26'
$ws.Range("E20").Value = '  +9.59%  '
$ws.Range("E21").Value = '  +1.74%  '
$ws.Range("E22").Value = '  +2.54%  '
$ws.Range("D23").Value = '66.61'
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("D24").Value = '241.21'
$ws.Range("E24").Value = '  +2.61%  '
$ws.Range("E25").Value = '  +4.17%  '
$ws.Range("E26").Value = '  +0.04%  '
$ws.Range("E27").Value = '  +5.11%  '
$ws.Range("D28").Value = '24.16'
$ws.Range("E28").Value = '  +4.72%  '
$ws.Range("D29").Value = '2.29'
$ws.Range("E29").Value = '  +9.96%  '
$ws.Range("E30").Value = '  +2.32%  '
$ws.Range("D31").Value = '34.49'
$ws.Range("E31").Value = '  +9.27%  '
$ws.Range("D32").Value = '160.18'
$ws.Range("E32").Value = '  +0.15%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("E34").Value = '  +4.02%  '
$ws.Range("D35").Value = '0.0743'
$ws.Range("E35").Value = '  +4.42%  '
$ws.Range("E36").Value = '  -1.46%  '
$ws.Range("E37").Value = '  +1.99%  '
$ws.Range("E38").Value = '  +2.84%  '
$ws.Range("D39").Value = '16.55'
$ws.Range("E39").Value = '  +6.69%  '
$ws.Range("E40").Value = '  +3.81%  '
$ws.Range("E41").Value = '  +2.62%  '
$ws.Range("E42").Value = '  +4.40%  '
$ws.Range("D43").Value = '2.059.55'
$ws.Range("E43").Value = '  -0.48%  '
$ws.Range("D44").Value = '19.31'
$ws.Range("E44").Value = '  +0.35%  '
$ws.Range("E45").Value = '  +2.82%  '
$ws.Range("D46").Value = '10.12'
$ws.Range("E46").Value = '  +2.50%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").Value = '2.04'
$ws.Range("E47").Value = '  +5.95%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").Value = '2.85'
$ws.Range("E48").Value = '  +3.54%  '
$ws.Range("D49").Value = '72.86'
$ws.Range("E49").Value = '  +7.90%  '
$ws.Range("E50").Value = '  +3.87%  '
$ws.Range("E51").Value = '  +3.02%  '
